$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.047038011799852
$ws.Range("D2").Value = 1.047035348639189
$ws.Range("E2").Value = 1.044322145195228
$ws.Range("F2").Value = 1.055215098588939
$ws.Range("I2").Value = 1.046144080025894
$ws.Range("J2").Value = 1.052089209779053
$ws.Range("K2").Value = 1.049799112658989
$ws.Range("L2").Value = 1.047093521497968
$ws.Range("M2").Value = 1.057956169872057
$ws.Range("N2").Value = 1.053583297133259
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.04865520034517
$ws.Range("D3").Value = 1.04790469875263
$ws.Range("E3").Value = 1.045742375168712
$ws.Range("F3").Value = 1.056945421500794
$ws.Range("I3").Value = 1.046662590946063
$ws.Range("J3").Value = 1.053351425694158
$ws.Range("K3").Value = 1.050480076317079
$ws.Range("L3").Value = 1.048323377233204
$ws.Range("M3").Value = 1.059497551721896
$ws.Range("N3").Value = 1.054847305539741
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.049698558764476
$ws.Range("D4").Value = 1.048465727194247
$ws.Range("E4").Value = 1.046658147257357
$ws.Range("F4").Value = 1.058062630481181
$ws.Range("I4").Value = 1.046995537729852
$ws.Range("J4").Value = 1.05416476110336
$ws.Range("K4").Value = 1.050918570783206
$ws.Range("L4").Value = 1.049115454070714
$ws.Range("M4").Value = 1.060492052253332
$ws.Range("N4").Value = 1.055661795978515
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.0501364652607
$ws.Range("D5").Value = 1.048701229202445
$ws.Range("E5").Value = 1.047042380957182
$ws.Range("F5").Value = 1.058531738350124
$ws.Range("I5").Value = 1.047134899229019
$ws.Range("J5").Value = 1.054505883730136
$ws.Range("K5").Value = 1.051102407435297
$ws.Range("L5").Value = 1.049447563280325
$ws.Range("M5").Value = 1.060909464426755
$ws.Range("N5").Value = 1.056003403038548
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.050209949783304
$ws.Range("D6").Value = 1.048740750378326
$ws.Range("E6").Value = 1.047106851448449
$ws.Range("F6").Value = 1.05861047079102
$ws.Range("I6").Value = 1.047158263039567
$ws.Range("J6").Value = 1.054563112916544
$ws.Range("K6").Value = 1.051133244853847
$ws.Range("L6").Value = 1.049503274603488
$ws.Range("M6").Value = 1.060979510505359
$ws.Range("N6").Value = 1.056060713496965
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.049704412908202
$ws.Range("D7").Value = 1.04846887536954
$ws.Range("E7").Value = 1.04666328436443
$ws.Range("F7").Value = 1.058068900931479
$ws.Range("I7").Value = 1.046997402272893
$ws.Range("J7").Value = 1.054169322342959
$ws.Range("K7").Value = 1.050921029201032
$ws.Range("L7").Value = 1.049119895172259
$ws.Range("M7").Value = 1.06049763237327
$ws.Range("N7").Value = 1.055666363695597
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.04758519236233
$ws.Range("D8").Value = 1.047329462767369
$ws.Range("E8").Value = 1.044802790812759
$ws.Range("F8").Value = 1.055800379977258
$ws.Range("I8").Value = 1.046319846992757
$ws.Range("J8").Value = 1.05251649277599
$ws.Range("K8").Value = 1.050029692270666
$ws.Range("L8").Value = 1.047509934435237
$ws.Range("M8").Value = 1.058477690447846
$ws.Range("N8").Value = 1.054011186921061
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.043826717215058
$ws.Range("D9").Value = 1.04530998308282
$ws.Range("E9").Value = 1.041499225124562
$ws.Range("F9").Value = 1.051783751334898
$ws.Range("I9").Value = 1.04510605922887
$ws.Range("J9").Value = 1.049577437837309
$ws.Range("K9").Value = 1.04844248662783
$ws.Range("L9").Value = 1.044643967464958
$ws.Range("M9").Value = 1.0548956857967
$ws.Range("N9").Value = 1.051067958187253
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.04130396734828
$ws.Range("D10").Value = 1.043955510947472
$ws.Range("E10").Value = 1.039279186295376
$ws.Range("F10").Value = 1.049092205285309
$ws.Range("I10").Value = 1.044283228629419
$ws.Range("J10").Value = 1.04759951213316
$ws.Range("K10").Value = 1.047372924162684
$ws.Range("L10").Value = 1.042713109491536
$ws.Range("M10").Value = 1.052491669979825
$ws.Range("N10").Value = 1.049087223601733
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.040207337987183
$ws.Range("D11").Value = 1.043367009512627
$ws.Range("E11").Value = 1.038313525581442
$ws.Range("F11").Value = 1.047923267170889
$ws.Range("I11").Value = 1.043923632899001
$ws.Range("J11").Value = 1.04673849048871
$ws.Range("K11").Value = 1.046907015485614
$ws.Range("L11").Value = 1.041872073648396
$ws.Range("M11").Value = 1.051446724786293
$ws.Range("N11").Value = 1.048224979207787
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.039799343789477
$ws.Range("D12").Value = 1.043148106655772
$ws.Range("E12").Value = 1.037954164987262
$ws.Range("F12").Value = 1.047488531442945
$ws.Range("I12").Value = 1.043789560771621
$ws.Range("J12").Value = 1.046417968988109
$ws.Range("K12").Value = 1.046733532616784
$ws.Range("L12").Value = 1.041558916282432
$ws.Range("M12").Value = 1.051057970322444
$ws.Range("N12").Value = 1.047904002529891
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.039886889944162
$ws.Range("D13").Value = 1.043195076058576
$ws.Range("E13").Value = 1.038031279667277
$ws.Range("F13").Value = 1.047581808448002
$ws.Range("I13").Value = 1.043818342519971
$ws.Range("J13").Value = 1.046486753822406
$ws.Range("K13").Value = 1.046770764545941
$ws.Range("L13").Value = 1.04162612417912
$ws.Range("M13").Value = 1.051141387572977
$ws.Range("N13").Value = 1.047972885046542
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.040173626565214
$ws.Range("D14").Value = 1.043348921235448
$ws.Range("E14").Value = 1.038283834501177
$ws.Range("F14").Value = 1.047887342896773
$ws.Range("I14").Value = 1.043912560734793
$ws.Range("J14").Value = 1.046712010435018
$ws.Range("K14").Value = 1.046892684013649
$ws.Range("L14").Value = 1.041846203530714
$ws.Range("M14").Value = 1.051414602864691
$ws.Range("N14").Value = 1.048198461549382
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.040350206863128
$ws.Range("D15").Value = 1.043443669440767
$ws.Range("E15").Value = 1.038439352482994
$ws.Range("F15").Value = 1.048075520675845
$ws.Range("I15").Value = 1.043970544957911
$ws.Range("J15").Value = 1.04685070533793
$ws.Range("K15").Value = 1.04696774630957
$ws.Range("L15").Value = 1.041981700665493
$ws.Range("M15").Value = 1.051582857816059
$ws.Range("N15").Value = 1.048337353414963
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.041376655792887
$ws.Range("D16").Value = 1.04399452508814
$ws.Range("E16").Value = 1.039343180703961
$ws.Range("F16").Value = 1.0491697090644
$ws.Range("I16").Value = 1.044307023697205
$ws.Range("J16").Value = 1.047656557928589
$ws.Range("K16").Value = 1.047403785898426
$ws.Range("L16").Value = 1.042768820562535
$ws.Range("M16").Value = 1.052560934134261
$ws.Range("N16").Value = 1.049144350408735
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.042019367149939
$ws.Range("D17").Value = 1.044339521110694
$ws.Range("E17").Value = 1.039908948115182
$ws.Range("F17").Value = 1.049855120583799
$ws.Range("I17").Value = 1.044517199116826
$ws.Range("J17").Value = 1.048160815750269
$ws.Range("K17").Value = 1.047676553358271
$ws.Range("L17").Value = 1.043261222036629
$ws.Range("M17").Value = 1.053173376127994
$ws.Range("N17").Value = 1.04964932433436
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.042393839739114
$ws.Range("D18").Value = 1.044540558296764
$ws.Range("E18").Value = 1.040238530111631
$ws.Range("F18").Value = 1.05025457483391
$ws.Range("I18").Value = 1.044639472426327
$ws.Range("J18").Value = 1.048454501108105
$ws.Range("K18").Value = 1.047835386066978
$ws.Range("L18").Value = 1.043547953702332
$ws.Range("M18").Value = 1.053530219232032
$ws.Range("N18").Value = 1.049943426759091
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.042521456168069
$ws.Range("D19").Value = 1.044609074194147
$ws.Range("E19").Value = 1.040350838165103
$ws.Range("F19").Value = 1.050390722154972
$ws.Range("I19").Value = 1.044681110622795
$ws.Range("J19").Value = 1.048554566054827
$ws.Range("K19").Value = 1.047889498644148
$ws.Range("L19").Value = 1.043645641164768
$ws.Range("M19").Value = 1.053651828858712
$ws.Range("N19").Value = 1.050043633809515
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.04195045284153
$ws.Range("D20").Value = 1.044302526323759
$ws.Range("E20").Value = 1.03984829023376
$ws.Range("F20").Value = 1.049781617156707
$ws.Range("I20").Value = 1.044494682264568
$ws.Range("J20").Value = 1.04810675923661
$ws.Range("K20").Value = 1.047647315747137
$ws.Range("L20").Value = 1.043208441537361
$ws.Range("M20").Value = 1.053107706725644
$ws.Range("N20").Value = 1.049595191054251
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.040089208056403
$ws.Range("D21").Value = 1.043303626193407
$ws.Range("E21").Value = 1.038209482045534
$ws.Range("F21").Value = 1.047797385663579
$ws.Range("I21").Value = 1.043884829730488
$ws.Range("J21").Value = 1.04664569740298
$ws.Range("K21").Value = 1.046856793510855
$ws.Range("L21").Value = 1.041781416716293
$ws.Range("M21").Value = 1.05133416488362
$ws.Range("N21").Value = 1.048132054345232
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.038915158568835
$ws.Range("D22").Value = 1.042673797599692
$ws.Range("E22").Value = 1.037175207442537
$ws.Range("F22").Value = 1.046546687308043
$ws.Range("I22").Value = 1.04349848239048
$ws.Range("J22").Value = 1.045723014584812
$ws.Range("K22").Value = 1.046357306807551
$ws.Range("L22").Value = 1.040879790095904
$ws.Range("M22").Value = 1.050215500899324
$ws.Range("N22").Value = 1.047208061211629
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.039537911350336
$ws.Range("D23").Value = 1.043007852387404
$ws.Range("E23").Value = 1.037723869763877
$ws.Range("F23").Value = 1.047210008925871
$ws.Range("I23").Value = 1.043703570097211
$ws.Range("J23").Value = 1.046212535183877
$ws.Range("K23").Value = 1.046622328732265
$ws.Range("L23").Value = 1.041358181088287
$ws.Range("M23").Value = 1.050808869344072
$ws.Range("N23").Value = 1.047698276986093
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.04198159351799
$ws.Range("D24").Value = 1.044319243272518
$ws.Range("E24").Value = 1.03987570021873
$ws.Range("F24").Value = 1.04981483122547
$ws.Range("I24").Value = 1.044504857644859
$ws.Range("J24").Value = 1.048131186406774
$ws.Range("K24").Value = 1.047660527793037
$ws.Range("L24").Value = 1.043232292249212
$ws.Range("M24").Value = 1.05313738108805
$ws.Range("N24").Value = 1.0496196529138
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.044801322664738
$ws.Range("D25").Value = 1.045833482137982
$ws.Range("E25").Value = 1.042356332812538
$ws.Range("F25").Value = 1.052824512150383
$ws.Range("I25").Value = 1.045422234644507
$ws.Range("J25").Value = 1.050340472803735
$ws.Range("K25").Value = 1.048854809339103
$ws.Range("L25").Value = 1.045388398963912
$ws.Range("M25").Value = 1.055824482161891
$ws.Range("N25").Value = 1.051832076750852
